$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly data refresh: insert a new "current week" record above the existing
# Albahaca price history (row 149), pushing all the older records down by one
# row (149 -> 150, 150 -> 151, ... 183 -> 184).
$ws.Rows.Item(149).Insert()

# Populate the newly inserted row 149 with this week's record.
$ws.Cells.Item(149, 1).Value = 4
$ws.Cells.Item(149, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(149, 3).Value = "Los Lagos"
$ws.Cells.Item(149, 4).Value = 45015
$ws.Cells.Item(149, 5).Value = 10
$ws.Cells.Item(149, 6).Value = 100112052
$ws.Cells.Item(149, 7).Value = "Albahaca"
$ws.Cells.Item(149, 8).Value = "Sin especificar"
$ws.Cells.Item(149, 9).Value = "Primera"
$ws.Cells.Item(149, 10).Value = 50
$ws.Cells.Item(149, 11).Value = 7000
$ws.Cells.Item(149, 12).Value = 7000
$ws.Cells.Item(149, 13).Value = 7000
$ws.Cells.Item(149, 14).Value = '$/docena de matas'
$ws.Cells.Item(149, 15).Value = "Región Metropolitana"
$ws.Cells.Item(149, 16).Value = 1167
$ws.Cells.Item(149, 17).Value = 6
$ws.Cells.Item(149, 18).Value = "Hortaliza"
